# Rebuild rows 2-10 of Sheet1 per the natmi LR-pairs edit (Thbs2-Notch3).
# Sending/Target clusters now form a full 3x3 cross of ECs/FAPs/sCs; rows 8-10 are new.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Thbs2", "Notch3", "ECs", 2, 0.6666666666666666, 0.8911683333333333, 2.673505, 0.02693425114262819, 0.02693425114262819, 3, 1, 7.457778999999999, 22.373337, 0.08029647035915141, 0.0802964703591514, 6.646136481798333, 59.815228336185, 0.002162725298519984, 0.002162725298519984),
    @("ECs", "Thbs2", "Notch3", "FAPs", 2, 0.6666666666666666, 0.8911683333333333, 2.673505, 0.02693425114262819, 0.02693425114262819, 3, 1, 2.732509, 8.197527, 0.02942039820764526, 0.02942039820764525, 2.435125491348333, 21.916129422135, 0.0007924163940408456, 0.0007924163940408454),
    @("ECs", "Thbs2", "Notch3", "sCs", 2, 0.6666666666666666, 0.8911683333333333, 2.673505, 0.02693425114262819, 0.02693425114262819, 3, 1, 82.68775466666666, 248.063264, 0.8902831314332034, 0.8902831314332033, 73.68870851336888, 663.19837662032, 0.02397910945006736, 0.02397910945006735),
    @("FAPs", "Thbs2", "Notch3", "ECs", 3, 1, 24.359699, 73.07909699999999, 0.7362360466408275, 0.7362360466408276, 3, 1, 7.457778999999999, 22.373337, 0.08029647035915141, 0.0802964703591514, 181.6692516485209, 1635.023264836689, 0.05911715589643402, 0.05911715589643402),
    @("FAPs", "Thbs2", "Notch3", "FAPs", 3, 1, 24.359699, 73.07909699999999, 0.7362360466408275, 0.7362360466408276, 3, 1, 2.732509, 8.197527, 0.02942039820764526, 0.02942039820764525, 66.56309675479099, 599.0678707931188, 0.02166035766699564, 0.02166035766699563),
    @("FAPs", "Thbs2", "Notch3", "sCs", 3, 1, 24.359699, 73.07909699999999, 0.7362360466408275, 0.7362360466408276, 3, 1, 82.68775466666666, 248.063264, 0.8902831314332034, 0.8902831314332033, 2014.248814665845, 18128.23933199261, 0.6554585330773979, 0.6554585330773979),
    @("sCs", "Thbs2", "Notch3", "ECs", 3, 1, 7.835938333333334, 23.507815, 0.2368297022165442, 0.2368297022165442, 3, 1, 7.457778999999999, 22.373337, 0.08029647035915141, 0.0802964703591514, 58.43869634762834, 525.948267128655, 0.0190165891641974, 0.0190165891641974),
    @("sCs", "Thbs2", "Notch3", "FAPs", 3, 1, 7.835938333333334, 23.507815, 0.2368297022165442, 0.2368297022165442, 3, 1, 2.732509, 8.197527, 0.02942039820764526, 0.02942039820764525, 21.41177201927833, 192.705948173505, 0.006967624146608779, 0.006967624146608777),
    @("sCs", "Thbs2", "Notch3", "sCs", 3, 1, 7.835938333333334, 23.507815, 0.2368297022165442, 0.2368297022165442, 3, 1, 82.68775466666666, 248.063264, 0.8902831314332034, 0.8902831314332033, 647.9361464897955, 5831.42531840816, 0.2108454889057381, 0.210845488905738)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $colLetter = [char](65 + $c)
        $ws.Range("$colLetter$r").Value = $rowVals[$c]
    }
}
